# Slide 6 ("SOURCES OF FINANCE") contains a graphicFrame (Shapes.Item(2))
# that hosts a table. The author changed the table's applied style from
# the custom "Table_0" style ({1DF289EE-91B6-45BA-A5DC-62FE9B2F4011}) to
# the built-in style {9CC2D14A-BD81-47D6-8730-2ADC76CD7405}
# (PowerPoint's "Medium Style 2 - Accent 1"), presumably picked from the
# Table Styles gallery on the Table Design ribbon tab.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(6)
$shape = $slide.Shapes.Item(2)
$table = $shape.Table

$table.ApplyStyle("{9CC2D14A-BD81-47D6-8730-2ADC76CD7405}")
